$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J (copy formatting from H1, which has the header style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2..32 (I = col 9, J = col 10)
$data = @{
    2  = @(9, 9)
    3  = @(7, 7)
    4  = @(8, 8)
    5  = @(8, 8)
    6  = @(7, 7)
    7  = @(1, 1)
    8  = @(9, 9)
    9  = @(6, 6)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(6, 6)
    13 = @(6, 6)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(8, 8)
    19 = @(8, 8)
    20 = @(6, 6)
    21 = @(9, 9)
    22 = @(8, 8)
    23 = @(6, 7)
    24 = @(7, 7)
    25 = @(7, 7)
    26 = @(7, 7)
    27 = @(6, 6)
    28 = @(6, 6)
    29 = @(6, 6)
    30 = @(5, 5)
    31 = @(5, 5)
    32 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
